# Updated symbol list on Mon Dec 26 06:43:20 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" column (D) with newly scraped quotes, and fixes a
# couple of "Bestin24h"/"Worstin24h" suffix labels in the "Volume(1h)"
# column (E) that coinranking appended/removed between scrapes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells are stored as text (General format, quote-prefixed numbers)
# in the source sheet, so every numeric-looking replacement is written
# with a leading apostrophe to keep it a text value instead of letting
# Excel coerce it into a floating point number.
$priceUpdates = [ordered]@{
    "D2"  = "243.56"
    "D3"  = "22.96"
    "D4"  = "5.416"
    "D5"  = "0.05945"
    "D6"  = "3.449"
    "D8"  = "0.8121"
    "D9"  = "0.9181"
    "D10" = "0.1408"
    "D11" = "0.07487"
    "D12" = "0.03282"
    "D13" = "0.03052"
    "D14" = "0.09348"
    "D15" = "3.860"
    "D16" = "0.001584"
    "D17" = "0.04678"
    "D18" = "0.0005936"
    "D19" = "0.006082"
    "D20" = "0.004994"
    "D21" = "0.0009829"
    "D22" = "0.0001099"
    "D25" = "0.3200"
    "D40" = "0.03942"
    "D41" = "0.006197"
    "D42" = "0.1073"
    "D44" = "0.008091"
    "D45" = "0.00005224"
    "D49" = "0.002275"
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $priceUpdates[$addr]
}

# Volume(1h) label tweaks.
$ws.Range("E22").Value = "21NitroExNTXBestin24h"
$ws.Range("E44").Value = "43LocalTradersLCT"
